$d = $word.ActiveDocument

# 1. Fix scenario title: "Scenario registrovanja korisnika" -> "Scenario učlanjenje u grupu"
$d.Content.Find.Execute("registrovanja korisnika", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "učlanjenje u grupu", 2)

# 2. Fix spelling: "zahtevom za pristupa i odobrenjem" -> "zahtevom za pristup i odobrenjem"
#    (only the first "pristupa" in that sentence should become "pristup")
$d.Content.Find.Execute("zahtevom za pristupa i odobrenjem", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "zahtevom za pristup i odobrenjem", 2)
